$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 794
$ws1.Range("F5").Value = 59
$ws1.Range("F6").Value = 86
$ws1.Range("F7").Value = 281
$ws1.Range("F8").Value = 3998
$ws1.Range("F10").Value = 4703
$ws1.Range("F11").Value = 519
$ws1.Range("F12").Value = 1186

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 794
$ws4.Range("F5").Value = 59
$ws4.Range("F6").Value = 86
$ws4.Range("F8").Value = 281
$ws4.Range("F9").Value = 3998
$ws4.Range("F11").Value = 4703
$ws4.Range("F12").Value = 519
$ws4.Range("F13").Value = 1186
